$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: new transaction number (new shared string)
$ws.Range("A18").Value = "OR.0026.0031"

# Row 19: repeat value "12105488" - copy an existing text-typed cell (A3)
# so the new cell keeps the same "text" storage type/style as the rest
# of the column instead of being coerced to a number.
$ws.Range("A3").Copy()
$ws.Range("A19").PasteSpecial()
